$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column to stay text (matches source inline-string cells; many
# values like "1.003" would otherwise be auto-parsed as numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Formula = "28.873.14"
$ws.Range("E2").Formula = "  +1.75%  "

# Row 3
$ws.Range("D3").Formula = "1.886.71"
$ws.Range("E3").Formula = "  +1.21%  "

# Row 4
$ws.Range("D4").Formula = "1.003"
$ws.Range("E4").Formula = "  +0.08%  "

# Row 5
$ws.Range("D5").Formula = "331.83"
$ws.Range("E5").Formula = "  -1.63%  "

# Row 6
$ws.Range("E6").Formula = "  +0.09%  "

# Row 7
$ws.Range("D7").Formula = "0.4621"
$ws.Range("E7").Formula = "  -1.69%  "

# Row 8
$ws.Range("D8").Formula = "0.4111"
$ws.Range("E8").Formula = "  +3.65%  "

# Row 9
$ws.Range("D9").Formula = "47.48"
$ws.Range("E9").Formula = "  -0.20%  "

# Row 10
$ws.Range("D10").Formula = "0.07970"
$ws.Range("E10").Formula = "  -0.46%  "

# Row 11
$ws.Range("D11").Formula = "0.9939"
$ws.Range("E11").Formula = "  -0.25%  "

# Row 12
$ws.Range("D12").Formula = "21.71"
$ws.Range("E12").Formula = "  -1.11%  "

# Row 13
$ws.Range("D13").Formula = "1.898.56"
$ws.Range("E13").Formula = "  +2.10%  "

# Row 14
$ws.Range("D14").Formula = "5.910"
$ws.Range("E14").Formula = "  -1.88%  "

# Row 15
$ws.Range("D15").Formula = "7.060"
$ws.Range("E15").Formula = "  -2.47%  "

# Row 16
$ws.Range("B16").Formula = "Litecoin"
$ws.Range("C16").Formula = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Formula = "89.24"
$ws.Range("E16").Formula = "  -1.22%  "

# Row 17
$ws.Range("B17").Formula = "BinanceUSD"
$ws.Range("C17").Formula = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Formula = "1.002"
$ws.Range("E17").Formula = "  -0.05%  "

# Row 18
$ws.Range("D18").Formula = "0.00001028"
$ws.Range("E18").Formula = "  -1.02%  "

# Row 19
$ws.Range("D19").Formula = "0.06567"
$ws.Range("E19").Formula = "  -1.09%  "

# Row 20
$ws.Range("D20").Formula = "17.47"
$ws.Range("E20").Formula = "  -0.23%  "

# Row 21
$ws.Range("E21").Formula = "  +0.09%  "

# Row 22
$ws.Range("D22").Formula = "28.907.45"
$ws.Range("E22").Formula = "  +1.77%  "

# Row 23
$ws.Range("D23").Formula = "5.386"
$ws.Range("E23").Formula = "  -1.46%  "

# Row 24
$ws.Range("D24").Formula = "11.25"
$ws.Range("E24").Formula = "  +1.94%  "

# Row 25
$ws.Range("D25").Formula = "2.216"
$ws.Range("E25").Formula = "  -2.33%  "

# Row 26
$ws.Range("D26").Formula = "2.118.55"
$ws.Range("E26").Formula = "  +1.66%  "

# Row 27
$ws.Range("D27").Formula = "157.60"
$ws.Range("E27").Formula = "  -1.88%  "

# Row 28
$ws.Range("D28").Formula = "19.68"
$ws.Range("E28").Formula = "  -0.16%  "

# Row 29
$ws.Range("D29").Formula = "2.123"
$ws.Range("E29").Formula = "  +0.68%  "

# Row 30
$ws.Range("D30").Formula = "5.417"
$ws.Range("E30").Formula = "  -0.93%  "

# Row 31
$ws.Range("D31").Formula = "117.81"
$ws.Range("E31").Formula = "  -1.48%  "

# Row 32
$ws.Range("D32").Formula = "0.9784"
$ws.Range("E32").Formula = "  +1.97%  "

# Row 33
$ws.Range("D33").Formula = "0.09362"
$ws.Range("E33").Formula = "  -1.54%  "

# Row 34
$ws.Range("D34").Formula = "1.414"
$ws.Range("E34").Formula = "  +2.88%  "

# Row 35
$ws.Range("D35").Formula = "3.606"
$ws.Range("E35").Formula = "  +0.31%  "

# Row 36
$ws.Range("D36").Formula = "5.278"
$ws.Range("E36").Formula = "  -1.37%  "

# Row 37
$ws.Range("D37").Formula = "0.06061"
$ws.Range("E37").Formula = "  -0.69%  "

# Row 38
$ws.Range("D38").Formula = "0.02235"
$ws.Range("E38").Formula = "  -0.84%  "

# Row 39
$ws.Range("D39").Formula = "8.296"
$ws.Range("E39").Formula = "  +0.25%  "

# Row 40
$ws.Range("D40").Formula = "1.174"
$ws.Range("E40").Formula = "  -0.40%  "

# Row 41
$ws.Range("E41").Formula = "  +0.08%  "

# Row 42
$ws.Range("D42").Formula = "0.5776"
$ws.Range("E42").Formula = "  -2.46%  "

# Row 43
$ws.Range("D43").Formula = "10.14"
$ws.Range("E43").Formula = "  -1.41%  "

# Row 44
$ws.Range("D44").Formula = "0.1816"
$ws.Range("E44").Formula = "  -3.09%  "

# Row 45
$ws.Range("D45").Formula = "1.266"
$ws.Range("E45").Formula = "  -0.70%  "

# Row 46
$ws.Range("D46").Formula = "2.283"
$ws.Range("E46").Formula = "  +10.93%  "

# Row 47
$ws.Range("D47").Formula = "0.5488"
$ws.Range("E47").Formula = "  -1.11%  "

# Row 48
$ws.Range("D48").Formula = "11.98"
$ws.Range("E48").Formula = "  -1.18%  "

# Row 49
$ws.Range("E49").Formula = "  -2.12%  "

# Row 50
$ws.Range("D50").Formula = "0.07012"
$ws.Range("E50").Formula = "  -4.25%  "

# Row 51
$ws.Range("B51").Formula = "Quant"
$ws.Range("C51").Formula = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Formula = "110.83"
$ws.Range("E51").Formula = "  -0.85%  "
